# Applies scheduled-runner profit recalculations to the Garuda_Profits sheets.
# Source data: Leve vendor/market price refresh across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 276.55554
$ws.Range("I4").Value = 44.142857
$ws.Range("J4").Value = 1090
$ws.Range("K4").Value = 44.142857
$ws.Range("L4").Value = 1090
$ws.Range("M4").Value = 69.85714300000001
$ws.Range("N4").Value = -1318
$ws.Range("H17").Value = 1724.6957
$ws.Range("J17").Value = 1851.4
$ws.Range("L17").Value = 5554.200000000001
$ws.Range("N17").Value = -5890.200000000001
$ws.Range("H19").Value = 2425.4285
$ws.Range("I19").Value = 2907.647
$ws.Range("J19").Value = 1680.1818
$ws.Range("K19").Value = 2907.647
$ws.Range("L19").Value = 1680.1818
$ws.Range("M19").Value = -2732.647
$ws.Range("N19").Value = -2030.1818
$ws.Range("H58").Value = 1529.3125
$ws.Range("J58").Value = 1911.3877
$ws.Range("L58").Value = 5734.1631
$ws.Range("N58").Value = -6034.1631
$ws.Range("H138").Value = 2969.5225
$ws.Range("J138").Value = 3455.0952
$ws.Range("L138").Value = 10365.2856
$ws.Range("N138").Value = -20645.2856

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 1800
$ws.Range("J8").Value = 1800
$ws.Range("L8").Value = 1800
$ws.Range("N8").Value = -2088
$ws.Range("H32").Value = 4179.65
$ws.Range("I32").Value = 3886.237
$ws.Range("J32").Value = 13666.667
$ws.Range("K32").Value = 3886.237
$ws.Range("L32").Value = 13666.667
$ws.Range("M32").Value = -3599.237
$ws.Range("N32").Value = -14240.667
$ws.Range("H61").Value = 1848.871
$ws.Range("I61").Value = 1490.9565
$ws.Range("J61").Value = 2877.875
$ws.Range("K61").Value = 1490.9565
$ws.Range("L61").Value = 2877.875
$ws.Range("M61").Value = -1278.9565
$ws.Range("N61").Value = -3301.875
$ws.Range("H132").Value = 8615.6875
$ws.Range("I132").Value = 9414.406999999999
$ws.Range("K132").Value = 28243.221
$ws.Range("M132").Value = -25713.221
$ws.Range("H136").Value = 1848.871
$ws.Range("I136").Value = 1490.9565
$ws.Range("J136").Value = 2877.875
$ws.Range("K136").Value = 4472.8695
$ws.Range("L136").Value = 8633.625
$ws.Range("M136").Value = -1922.8695
$ws.Range("N136").Value = -13733.625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1886.3429
$ws.Range("I107").Value = 1487.625
$ws.Range("J107").Value = 2756.2727
$ws.Range("K107").Value = 1487.625
$ws.Range("L107").Value = 2756.2727
$ws.Range("M107").Value = 432.375
$ws.Range("N107").Value = -6596.2727

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 784.625
$ws.Range("I10").Value = 716.1667
$ws.Range("J10").Value = 990
$ws.Range("K10").Value = 716.1667
$ws.Range("L10").Value = 990
$ws.Range("M10").Value = -577.1667
$ws.Range("N10").Value = -1268
$ws.Range("H58").Value = 1167.5682
$ws.Range("I58").Value = 1329.4166
$ws.Range("J58").Value = 439.25
$ws.Range("K58").Value = 1329.4166
$ws.Range("L58").Value = 439.25
$ws.Range("M58").Value = -1126.4166
$ws.Range("N58").Value = -845.25
$ws.Range("H132").Value = 5436437.5
$ws.Range("I132").Value = 1243.3684
$ws.Range("J132").Value = 31253610
$ws.Range("K132").Value = 3730.1052
$ws.Range("L132").Value = 93760830
$ws.Range("M132").Value = -1200.1052
$ws.Range("N132").Value = -93765890
$ws.Range("H136").Value = 1167.5682
$ws.Range("I136").Value = 1329.4166
$ws.Range("J136").Value = 439.25
$ws.Range("K136").Value = 3988.2498
$ws.Range("L136").Value = 1317.75
$ws.Range("M136").Value = -1438.2498
$ws.Range("N136").Value = -6417.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1272.909
$ws.Range("I34").Value = 399.66666
$ws.Range("J34").Value = 1600.375
$ws.Range("K34").Value = 1198.99998
$ws.Range("L34").Value = 4801.125
$ws.Range("M34").Value = -1114.99998
$ws.Range("N34").Value = -4969.125
$ws.Range("H38").Value = 42.4
$ws.Range("I38").Value = 25.5
$ws.Range("J38").Value = 53.666668
$ws.Range("K38").Value = 76.5
$ws.Range("L38").Value = 161.000004
$ws.Range("M38").Value = 270.5
$ws.Range("N38").Value = -855.000004
$ws.Range("H39").Value = 3333
$ws.Range("J39").Value = 3333
$ws.Range("L39").Value = 9999
$ws.Range("N39").Value = -10587
$ws.Range("H42").Value = 2690
$ws.Range("I42").Value = 3000
$ws.Range("J42").Value = 2612.5
$ws.Range("K42").Value = 9000
$ws.Range("L42").Value = 7837.5
$ws.Range("M42").Value = -8466
$ws.Range("N42").Value = -8905.5
$ws.Range("H55").Value = 2596
$ws.Range("J55").Value = 2596
$ws.Range("L55").Value = 7788
$ws.Range("N55").Value = -8142
$ws.Range("H139").Value = 33335432
$ws.Range("I139").Value = 38463268
$ws.Range("K139").Value = 115389804
$ws.Range("M139").Value = -115384664

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5865.6787
$ws.Range("I132").Value = 6049.64
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 18148.92
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -15618.92
$ws.Range("N132").Value = -18057.9995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4626.8335
$ws.Range("I122").Value = 6375.8887
$ws.Range("J122").Value = 2877.7778
$ws.Range("K122").Value = 19127.6661
$ws.Range("L122").Value = 8633.3334
$ws.Range("M122").Value = -16677.6661
$ws.Range("N122").Value = -13533.3334
$ws.Range("H132").Value = 6637.6904
$ws.Range("I132").Value = 8621.679
$ws.Range("J132").Value = 2669.7144
$ws.Range("K132").Value = 25865.037
$ws.Range("L132").Value = 8009.1432
$ws.Range("M132").Value = -23335.037
$ws.Range("N132").Value = -13069.1432
$ws.Range("H133").Value = 25689.285
$ws.Range("J133").Value = 25689.285
$ws.Range("L133").Value = 25689.285
$ws.Range("N133").Value = -30749.285

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2148
$ws.Range("I122").Value = 2132
$ws.Range("J122").Value = 2191.077
$ws.Range("K122").Value = 6396
$ws.Range("L122").Value = 6573.231000000001
$ws.Range("M122").Value = -3946
$ws.Range("N122").Value = -11473.231
$ws.Range("H136").Value = 1849.0638
$ws.Range("I136").Value = 1732.1714
$ws.Range("J136").Value = 2190
$ws.Range("K136").Value = 5196.5142
$ws.Range("L136").Value = 6570
$ws.Range("M136").Value = -2646.5142
$ws.Range("N136").Value = -11670

